# Updates cryptos list values per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the assigned value to remain a text cell (matches the workbook's
    # original inlineStr cells) rather than being auto-coerced to a number by
    # Excel's normal type inference, then drop back to the Normal style so no
    # stray number-format/quote-prefix is left on the cell.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "63.183.62"
$ws.Range("E2").Value = "  +0.05%  "

Set-TextValue $ws.Range("D3") "2.612.53"
$ws.Range("E3").Value = "  -1.48%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue $ws.Range("D5") "606.85"
$ws.Range("E5").Value = "  +2.46%  "

Set-TextValue $ws.Range("D6") "145.74"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("E7").Value = "  -0.16%  "

Set-TextValue $ws.Range("D8") "0.586"
$ws.Range("E8").Value = "  -0.04%  "

Set-TextValue $ws.Range("D9") "2.609.28"
$ws.Range("E9").Value = "  -1.89%  "

$ws.Range("E10").Value = "  +1.37%  "

$ws.Range("E11").Value = "  -2.91%  "

Set-TextValue $ws.Range("D12") "0.372"
$ws.Range("E12").Value = "  +5.01%  "

$ws.Range("E13").Value = "  -0.20%  "

Set-TextValue $ws.Range("D14") "27.15"
$ws.Range("E14").Value = "  -0.62%  "

Set-TextValue $ws.Range("D15") "3.082.40"
$ws.Range("E15").Value = "  -1.58%  "

Set-TextValue $ws.Range("D16") "63.025.81"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("E17").Value = "  +1.81%  "

Set-TextValue $ws.Range("D18") "2.626.54"
$ws.Range("E18").Value = "  -0.87%  "

Set-TextValue $ws.Range("D19") "11.48"
$ws.Range("E19").Value = "  +0.65%  "

Set-TextValue $ws.Range("D20") "4.52"
$ws.Range("E20").Value = "  +2.97%  "

Set-TextValue $ws.Range("D21") "341.95"
$ws.Range("E21").Value = "  +1.12%  "

Set-TextValue $ws.Range("D22") "6.85"
$ws.Range("E22").Value = "  +1.26%  "

Set-TextValue $ws.Range("D24") "5.71"
$ws.Range("E24").Value = "  -1.15%  "

Set-TextValue $ws.Range("D25") "66.19"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("E26").Value = "  +0.80%  "

Set-TextValue $ws.Range("D27") "1.60"
$ws.Range("E27").Value = "  +5.31%  "

Set-TextValue $ws.Range("D28") "9.02"
$ws.Range("E28").Value = "  +7.34%  "

Set-TextValue $ws.Range("D29") "554.48"
$ws.Range("E29").Value = "  +3.76%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D30") "0.161"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D31") "1.00"
$ws.Range("E31").Value = "  -0.10%  "

Set-TextValue $ws.Range("D32") "7.76"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("E33").Value = "  +2.26%  "

Set-TextValue $ws.Range("D34") "0.0₃0841"
$ws.Range("E34").Value = "  +4.49%  "

$ws.Range("E35").Value = "  -4.25%  "

Set-TextValue $ws.Range("D36") "5.17"
$ws.Range("E36").Value = "  +2.55%  "

Set-TextValue $ws.Range("D37") "168.19"
$ws.Range("E37").Value = "  -2.82%  "

Set-TextValue $ws.Range("D38") "0.998"
$ws.Range("E38").Value = "  -0.26%  "

Set-TextValue $ws.Range("D39") "0.402"
$ws.Range("E39").Value = "  -1.02%  "

Set-TextValue $ws.Range("D40") "1.93"
$ws.Range("E40").Value = "  +5.84%  "

$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("E42").Value = "  +0.09%  "

Set-TextValue $ws.Range("D43") "165.30"
$ws.Range("E43").Value = "  -4.05%  "

Set-TextValue $ws.Range("D44") "39.60"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("E45").Value = "  -0.26%  "

Set-TextValue $ws.Range("D46") "21.72"
$ws.Range("E46").Value = "  -1.26%  "

Set-TextValue $ws.Range("D47") "0.0562"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("E49").Value = "  +2.57%  "

$ws.Range("E50").Value = "  -0.43%  "

Set-TextValue $ws.Range("D51") "1.92"
$ws.Range("E51").Value = "  +12.89%  "
